# Applies the recorded change: the observation rows got re-sorted in the
# source system, which (because the sheet was re-exported with the new
# ordering) shows up as a handful of row-groups trading their
# per-observation data while the shared location/metadata columns stay put.
#
# Row groups affected: (3 4), (6 7), (11 12 13) [3-cycle], (14 15).
#
# NOTE: this interpreter's PowerShell does not bind named ("-Param value")
# function arguments, so helper functions below use positional params.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowFields {
    param($Row, $Fields)
    foreach ($col in $Fields.Keys) {
        $cellRef = "$col$Row"
        $val = $Fields[$col]
        if ($null -eq $val) {
            $ws.Range($cellRef).Value = ""
        } else {
            $ws.Range($cellRef).Value = $val
        }
    }
}

# ---- Row 3 <- old Row 4 data; Row 4 <- old Row 3 data ----
Set-RowFields 3 @{
    A  = 131047025
    B  = 89193
    E  = 510
    F  = "Doftskinn"
    G  = "Cystostereum murrayi"
    H  = "(Berk. & M.A.Curtis.) Pouzar"
    Q  = 402314
    R  = 6818423
    Z  = "16:05"
    AB = "16:05"
}

Set-RowFields 4 @{
    A  = 131046847
    B  = 79243
    E  = 6425
    F  = "Garnlav"
    G  = "Alectoria sarmentosa"
    H  = "(Ach.) Ach."
    Q  = 402380
    R  = 6818405
    Z  = "17:01"
    AB = "17:01"
}

# ---- Row 6 <- old Row 7 data; Row 7 <- old Row 6 data ----
Set-RowFields 6 @{
    A  = 131046755
    B  = 57881
    E  = 100049
    F  = "Spillkråka"
    G  = "Dryocopus martius"
    H  = "(Linnaeus, 1758)"
    M  = "färska spår"
    Q  = 402424
    R  = 6818357
    Z  = "16:56"
    AB = "16:56"
}

Set-RowFields 7 @{
    A  = 131046733
    B  = 91808
    E  = 1202
    F  = "Ullticka"
    G  = "Phellinidium ferrugineofuscum"
    H  = "(P.Karst.) Fiasson & Niemelä"
    M  = $null
    Q  = 402493
    R  = 6818443
    Z  = "16:43"
    AB = "16:43"
}

# ---- 3-cycle: Row 11 <- old Row 13; Row 12 <- old Row 11; Row 13 <- old Row 12 ----
Set-RowFields 11 @{
    A  = 131046763
    B  = 92267
    D  = "VU"
    E  = 1209
    F  = "Rynkskinn"
    G  = "Hermanssonia centrifuga"
    H  = "(P. Karst.) Zmitr."
    M  = $null
    Q  = 402378
    R  = 6818392
    Z  = "17:02"
    AB = "17:02"
}

Set-RowFields 12 @{
    A  = 131046735
    M  = "nyligen använt bo"
    Q  = 402448
    R  = 6818295
    Z  = "16:54"
    AB = "16:54"
    AC = $null
}

Set-RowFields 13 @{
    A  = 131046788
    B  = 57884
    D  = "NT"
    E  = 100109
    F  = "Tretåig hackspett"
    G  = "Picoides tridactylus"
    H  = "(Linnaeus, 1758)"
    M  = "färska spår"
    Q  = 402473
    R  = 6818425
    Z  = "16:47"
    AB = "16:47"
    AC = "Färska ringhack (gran)"
}

# ---- Row 14 <- old Row 15 data; Row 15 <- old Row 14 data ----
Set-RowFields 14 @{
    A  = 131046808
    Q  = 402323
    R  = 6818416
    Z  = "16:06"
    AB = "16:06"
}

Set-RowFields 15 @{
    A  = 131046811
    Q  = 402450
    R  = 6818298
    Z  = "16:54"
    AB = "16:54"
}
